$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 8.324917666666666
$ws.Range("H2").Value = 24.974753
$ws.Range("I2").Value = 0.8193616330571973
$ws.Range("J2").Value = 0.8193616330571972
$ws.Range("M2").Value = 38.55267666666666
$ws.Range("N2").Value = 115.65803
$ws.Range("O2").Value = 0.5758151725879548
$ws.Range("P2").Value = 0.5758151725879548
$ws.Range("Q2").Value = 320.9478590796211
$ws.Range("R2").Value = 2888.53073171659
$ws.Range("S2").Value = 0.4718008601507785
$ws.Range("T2").Value = 0.4718008601507785
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 8.324917666666666
$ws.Range("H3").Value = 24.974753
$ws.Range("I3").Value = 0.8193616330571973
$ws.Range("J3").Value = 0.8193616330571972
$ws.Range("O3").Value = 0.08021535714867321
$ws.Range("P3").Value = 0.08021535714867323
$ws.Range("Q3").Value = 44.71043551434267
$ws.Range("R3").Value = 402.3939196290841
$ws.Range("S3").Value = 0.06572538602960322
$ws.Range("T3").Value = 0.06572538602960322
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 8.324917666666666
$ws.Range("H4").Value = 24.974753
$ws.Range("I4").Value = 0.8193616330571973
$ws.Range("J4").Value = 0.8193616330571972
$ws.Range("M4").Value = 23.02986166666667
$ws.Range("N4").Value = 69.089585
$ws.Range("O4").Value = 0.3439694702633719
$ws.Range("P4").Value = 0.3439694702633719
$ws.Range("Q4").Value = 191.7217022497227
$ws.Range("R4").Value = 1725.495320247505
$ws.Range("S4").Value = 0.2818353868768155
$ws.Range("T4").Value = 0.2818353868768155
$ws.Range("G5").Value = 1.835330666666667
$ws.Range("H5").Value = 5.505992
$ws.Range("I5").Value = 0.1806383669428028
$ws.Range("J5").Value = 0.1806383669428027
$ws.Range("M5").Value = 38.55267666666666
$ws.Range("N5").Value = 115.65803
$ws.Range("O5").Value = 0.5758151725879548
$ws.Range("P5").Value = 0.5758151725879548
$ws.Range("Q5").Value = 70.75690976841777
$ws.Range("R5").Value = 636.81218791576
$ws.Range("S5").Value = 0.1040143124371763
$ws.Range("T5").Value = 0.1040143124371763
$ws.Range("G6").Value = 1.835330666666667
$ws.Range("H6").Value = 5.505992
$ws.Range("I6").Value = 0.1806383669428028
$ws.Range("J6").Value = 0.1806383669428027
$ws.Range("O6").Value = 0.08021535714867321
$ws.Range("P6").Value = 0.08021535714867323
$ws.Range("Q6").Value = 9.856966363530667
$ws.Range("R6").Value = 88.71269727177601
$ws.Range("S6").Value = 0.01448997111907001
$ws.Range("T6").Value = 0.01448997111907001
$ws.Range("G7").Value = 1.835330666666667
$ws.Range("H7").Value = 5.505992
$ws.Range("I7").Value = 0.1806383669428028
$ws.Range("J7").Value = 0.1806383669428027
$ws.Range("M7").Value = 23.02986166666667
$ws.Range("N7").Value = 69.089585
$ws.Range("O7").Value = 0.3439694702633719
$ws.Range("P7").Value = 0.3439694702633719
$ws.Range("Q7").Value = 42.26741136592445
$ws.Range("R7").Value = 380.40670229332
$ws.Range("S7").Value = 0.06213408338655646
$ws.Range("T7").Value = 0.06213408338655646
